# Update hainan aqua chart2
# - Rename sheet 2 ("左2-近年淡水养殖面积情况") -> "左2-近年水产品产量及构成"
# - Replace its single-series content (year, 淡水养殖) with a 3-series
#   production table (year, 总产量, 海水产品, 淡水产品)
# - Make sheet 2 the active tab (was sheet 1)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- sheet1: restyle B1:D6 from the (soon to be orphaned/removed) duplicate
#     red-font style to the plain data style already used elsewhere (s=10).
#     Source the format from sheet3!B2, which keeps that style untouched.
$fmtSrc10 = $ws3.Range("B2")
$fmtSrc10.Copy()
$ws1.Range("B1:D6").PasteSpecial(-4122)  # xlPasteFormats

# --- sheet2: rename and rebuild content ---
$ws2.Name = "左2-近年水产品产量及构成"

# Reuse the existing header/data style (s=13, centered+wrap, used on sheet1!A2 etc.)
$fmtSrc13 = $ws1.Range("A2")
$fmtSrc13.Copy()
$ws2.Range("A1:D6").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("A1").Value = ""
$ws2.Range("B1").Value = "总产量(万吨)"
$ws2.Range("C1").Value = "海水产品万吨)`n"
$ws2.Range("D1").Value = "淡水产品(万吨)`n"

$ws2.Range("A2").Value = 2018
$ws2.Range("B2").Value = 175.8188
$ws2.Range("C2").Value = 137.8071
$ws2.Range("D2").Value = 38.0117

$ws2.Range("A3").Value = 2019
$ws2.Range("B3").Value = 172.1571
$ws2.Range("C3").Value = 135.0103
$ws2.Range("D3").Value = 37.1468

$ws2.Range("A4").Value = 2020
$ws2.Range("B4").Value = 166.7878
$ws2.Range("C4").Value = 130.2591
$ws2.Range("D4").Value = 36.5287

$ws2.Range("A5").Value = 2021
$ws2.Range("B5").Value = 164.0918
$ws2.Range("C5").Value = 127.3762
$ws2.Range("D5").Value = 36.7156

$ws2.Range("A6").Value = 2022
$ws2.Range("B6").Value = 170.311
$ws2.Range("C6").Value = 128.0529
$ws2.Range("D6").Value = 42.2581

$ws2.Rows.Item(1).RowHeight = 39.75

# --- activate sheet2 and move its selection ---
$ws2.Activate()
$ws2.Range("G20").Select()
